$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 855-870 (only changed cells) ---

# Row 855
$ws.Cells.Item(855, 4).Value = 44890
$ws.Cells.Item(855, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(855, 11).Value = 'Florida King'
$ws.Cells.Item(855, 13).Value = 220
$ws.Cells.Item(855, 14).Value = 10000
$ws.Cells.Item(855, 15).Value = 10000
$ws.Cells.Item(855, 16).Value = 10000
$ws.Cells.Item(855, 19).Value = 1250

# Row 856
$ws.Cells.Item(856, 4).Value = 44890
$ws.Cells.Item(856, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(856, 11).Value = 'Florida King'
$ws.Cells.Item(856, 12).Value = 'Primera'
$ws.Cells.Item(856, 13).Value = 180
$ws.Cells.Item(856, 14).Value = 8000
$ws.Cells.Item(856, 15).Value = 8000
$ws.Cells.Item(856, 16).Value = 8000
$ws.Cells.Item(856, 17).Value = '$/bandeja 8 kilos empedrada'
$ws.Cells.Item(856, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(856, 19).Value = 1000
$ws.Cells.Item(856, 20).Value = 8

# Row 857
$ws.Cells.Item(857, 4).Value = 44890
$ws.Cells.Item(857, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(857, 11).Value = 'Florida King'
$ws.Cells.Item(857, 12).Value = 'Segunda'
$ws.Cells.Item(857, 13).Value = 200
$ws.Cells.Item(857, 14).Value = 6000
$ws.Cells.Item(857, 15).Value = 6000
$ws.Cells.Item(857, 16).Value = 6000
$ws.Cells.Item(857, 17).Value = '$/bandeja 8 kilos empedrada'
$ws.Cells.Item(857, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(857, 19).Value = 750
$ws.Cells.Item(857, 20).Value = 8

# Row 858
$ws.Cells.Item(858, 12).Value = 'Especial'
$ws.Cells.Item(858, 13).Value = 120
$ws.Cells.Item(858, 14).Value = 8000
$ws.Cells.Item(858, 15).Value = 8000
$ws.Cells.Item(858, 16).Value = 8000
$ws.Cells.Item(858, 19).Value = 1000

# Row 859
$ws.Cells.Item(859, 12).Value = 'Especial'
$ws.Cells.Item(859, 13).Value = 5
$ws.Cells.Item(859, 14).Value = 270000
$ws.Cells.Item(859, 15).Value = 270000
$ws.Cells.Item(859, 16).Value = 270000
$ws.Cells.Item(859, 19).Value = 675

# Row 860
$ws.Cells.Item(860, 12).Value = 'Especial'
$ws.Cells.Item(860, 13).Value = 80
$ws.Cells.Item(860, 14).Value = 15000
$ws.Cells.Item(860, 15).Value = 15000
$ws.Cells.Item(860, 16).Value = 15000
$ws.Cells.Item(860, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(860, 19).Value = 938
$ws.Cells.Item(860, 20).Value = 16

# Row 861
$ws.Cells.Item(861, 13).Value = 250
$ws.Cells.Item(861, 14).Value = 6000
$ws.Cells.Item(861, 15).Value = 6000
$ws.Cells.Item(861, 16).Value = 6000
$ws.Cells.Item(861, 17).Value = '$/bandeja 8 kilos empedrada'
$ws.Cells.Item(861, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(861, 19).Value = 750
$ws.Cells.Item(861, 20).Value = 8

# Row 862
$ws.Cells.Item(862, 12).Value = 'Primera'
$ws.Cells.Item(862, 13).Value = 8
$ws.Cells.Item(862, 14).Value = 250000
$ws.Cells.Item(862, 15).Value = 250000
$ws.Cells.Item(862, 16).Value = 250000
$ws.Cells.Item(862, 17).Value = '$/bins (400 kilos)'
$ws.Cells.Item(862, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(862, 20).Value = 400

# Row 863
$ws.Cells.Item(863, 12).Value = 'Primera'
$ws.Cells.Item(863, 13).Value = 8
$ws.Cells.Item(863, 14).Value = 250000
$ws.Cells.Item(863, 15).Value = 250000
$ws.Cells.Item(863, 16).Value = 250000
$ws.Cells.Item(863, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(863, 19).Value = 625

# Row 864
$ws.Cells.Item(864, 12).Value = 'Primera'
$ws.Cells.Item(864, 13).Value = 120
$ws.Cells.Item(864, 14).Value = 13000
$ws.Cells.Item(864, 15).Value = 13000
$ws.Cells.Item(864, 16).Value = 13000
$ws.Cells.Item(864, 19).Value = 812

# Row 865
$ws.Cells.Item(865, 4).Value = 44225
$ws.Cells.Item(865, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(865, 12).Value = 'Segunda'
$ws.Cells.Item(865, 13).Value = 190
$ws.Cells.Item(865, 14).Value = 5000
$ws.Cells.Item(865, 15).Value = 5000
$ws.Cells.Item(865, 16).Value = 5000
$ws.Cells.Item(865, 17).Value = '$/bandeja 8 kilos empedrada'
$ws.Cells.Item(865, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(865, 19).Value = 625
$ws.Cells.Item(865, 20).Value = 8

# Row 866
$ws.Cells.Item(866, 4).Value = 44225
$ws.Cells.Item(866, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(866, 12).Value = 'Segunda'
$ws.Cells.Item(866, 13).Value = 12
$ws.Cells.Item(866, 14).Value = 230000
$ws.Cells.Item(866, 15).Value = 230000
$ws.Cells.Item(866, 16).Value = 230000
$ws.Cells.Item(866, 17).Value = '$/bins (400 kilos)'
$ws.Cells.Item(866, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(866, 19).Value = 575
$ws.Cells.Item(866, 20).Value = 400

# Row 867
$ws.Cells.Item(867, 4).Value = 44225
$ws.Cells.Item(867, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(867, 12).Value = 'Segunda'
$ws.Cells.Item(867, 13).Value = 140
$ws.Cells.Item(867, 14).Value = 12000
$ws.Cells.Item(867, 15).Value = 12000
$ws.Cells.Item(867, 16).Value = 12000
$ws.Cells.Item(867, 19).Value = 750

# Row 868
$ws.Cells.Item(868, 4).Value = 44225
$ws.Cells.Item(868, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(868, 12).Value = 'Especial'
$ws.Cells.Item(868, 13).Value = 100
$ws.Cells.Item(868, 14).Value = 14000
$ws.Cells.Item(868, 15).Value = 14000
$ws.Cells.Item(868, 16).Value = 14000
$ws.Cells.Item(868, 19).Value = 875

# Row 869
$ws.Cells.Item(869, 4).Value = 44225
$ws.Cells.Item(869, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(869, 12).Value = 'Especial'
$ws.Cells.Item(869, 13).Value = 155
$ws.Cells.Item(869, 14).Value = 13000
$ws.Cells.Item(869, 15).Value = 13000
$ws.Cells.Item(869, 16).Value = 13000
$ws.Cells.Item(869, 19).Value = 812

# Row 870
$ws.Cells.Item(870, 12).Value = 'Extra (doble especial)'
$ws.Cells.Item(870, 13).Value = 60
$ws.Cells.Item(870, 14).Value = 16000
$ws.Cells.Item(870, 15).Value = 16000
$ws.Cells.Item(870, 16).Value = 16000
$ws.Cells.Item(870, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(870, 19).Value = 1000

# --- Append new rows 871-873 ---

# Row 871
$ws.Cells.Item(871, 1).Value = 9
$ws.Cells.Item(871, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(871, 3).Value = 'Metropolitana'
$ws.Cells.Item(871, 4).Value = 44236
$ws.Cells.Item(871, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(871, 5).Value = 13
$ws.Cells.Item(871, 6).Value = 'Fruta'
$ws.Cells.Item(871, 7).Value = 100103
$ws.Cells.Item(871, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(871, 9).Value = 100103004
$ws.Cells.Item(871, 10).Value = 'Durazno'
$ws.Cells.Item(871, 11).Value = 'Carson'
$ws.Cells.Item(871, 12).Value = 'Primera'
$ws.Cells.Item(871, 13).Value = 125
$ws.Cells.Item(871, 14).Value = 10000
$ws.Cells.Item(871, 15).Value = 10000
$ws.Cells.Item(871, 16).Value = 10000
$ws.Cells.Item(871, 17).Value = '$/caja 16 kilos empedrada'
$ws.Cells.Item(871, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(871, 19).Value = 625
$ws.Cells.Item(871, 20).Value = 16

# Row 872
$ws.Cells.Item(872, 1).Value = 9
$ws.Cells.Item(872, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(872, 3).Value = 'Metropolitana'
$ws.Cells.Item(872, 4).Value = 44236
$ws.Cells.Item(872, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(872, 5).Value = 13
$ws.Cells.Item(872, 6).Value = 'Fruta'
$ws.Cells.Item(872, 7).Value = 100103
$ws.Cells.Item(872, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(872, 9).Value = 100103004
$ws.Cells.Item(872, 10).Value = 'Durazno'
$ws.Cells.Item(872, 11).Value = 'Carson'
$ws.Cells.Item(872, 12).Value = 'Primera'
$ws.Cells.Item(872, 13).Value = 150
$ws.Cells.Item(872, 14).Value = 9000
$ws.Cells.Item(872, 15).Value = 9000
$ws.Cells.Item(872, 16).Value = 9000
$ws.Cells.Item(872, 17).Value = '$/caja 16 kilos granel'
$ws.Cells.Item(872, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(872, 19).Value = 562
$ws.Cells.Item(872, 20).Value = 16

# Row 873
$ws.Cells.Item(873, 1).Value = 9
$ws.Cells.Item(873, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(873, 3).Value = 'Metropolitana'
$ws.Cells.Item(873, 4).Value = 44236
$ws.Cells.Item(873, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(873, 5).Value = 13
$ws.Cells.Item(873, 6).Value = 'Fruta'
$ws.Cells.Item(873, 7).Value = 100103
$ws.Cells.Item(873, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(873, 9).Value = 100103004
$ws.Cells.Item(873, 10).Value = 'Durazno'
$ws.Cells.Item(873, 11).Value = 'Carson'
$ws.Cells.Item(873, 12).Value = 'Segunda'
$ws.Cells.Item(873, 13).Value = 135
$ws.Cells.Item(873, 14).Value = 8000
$ws.Cells.Item(873, 15).Value = 8000
$ws.Cells.Item(873, 16).Value = 8000
$ws.Cells.Item(873, 17).Value = '$/caja 16 kilos granel'
$ws.Cells.Item(873, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(873, 19).Value = 500
$ws.Cells.Item(873, 20).Value = 16

